$d = $word.ActiveDocument

# The body currently ends with:
#   ... <empty "Geenafstand" paragraph> <wrapper table containing a nested
#   table that holds the signature image> <trailing empty paragraph>
# The edit removes the wrapper table (and its nested table/image) together
# with the trailing empty paragraph, and replaces them with a single new
# "Geenafstand" paragraph that directly holds the (re-uploaded) image.

# 1) Drop the outer signature table (this also removes the nested table it
#    contains, since that lives inside one of its cells).
if ($d.Tables.Count -ge 1) {
    $d.Tables.Item(1).Delete()
}

# 2) The document body now ends with the old trailing empty paragraph mark.
#    Re-fetch Content fresh (indexing Paragraphs right after a table delete
#    is unreliable) and target that paragraph mark directly by range offset.
$content = $d.Content
$tail = $d.Range($content.End - 1, $content.End)

# 3) Replace that paragraph with the new one: "Geenafstand" style carrying
#    a run with the inline picture (new docPr/cNvPr ids, alt-text descr,
#    and anchorId/editId, same embedded image relationship rId4).
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:pPr><w:pStyle w:val="Geenafstand"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Aptos"/><w:noProof/><w:sz w:val="2"/><w:szCs w:val="2"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="28BAD001" wp14:editId="6FAC9EDC"><wp:extent cx="1437005" cy="451485"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1356270451" name="Afbeelding 1401568752" descr="Afbeelding met Graphics, grafische vormgeving, creativiteit&#10;&#10;Door AI gegenereerde inhoud is mogelijk onjuist."/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1356270451" name="Afbeelding 1401568752" descr="Afbeelding met Graphics, grafische vormgeving, creativiteit&#10;&#10;Door AI gegenereerde inhoud is mogelijk onjuist."/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="1437005" cy="451485"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$tail.InsertXML($xml)
